$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 '68.936.71'
Set-TextValue 2 5 '  -2.50%  '
Set-TextValue 3 4 '3.508.14'
Set-TextValue 3 5 '  -3.28%  '
Set-TextValue 4 4 '1.00'
Set-TextValue 4 5 '  -0.03%  '
Set-TextValue 5 4 '585.29'
Set-TextValue 5 5 '  +0.25%  '
Set-TextValue 6 4 '168.27'
Set-TextValue 6 5 '  -4.32%  '
Set-TextValue 7 4 '0.609'
Set-TextValue 7 5 '  -3.92%  '
Set-TextValue 8 4 '3.501.22'
Set-TextValue 8 5 '  -3.28%  '
Set-TextValue 9 4 '1.00'
Set-TextValue 9 5 '  +0.03%  '
Set-TextValue 10 4 '0.187'
Set-TextValue 10 5 '  -4.28%  '
Set-TextValue 11 4 '6.70'
Set-TextValue 11 5 '  -1.10%  '
Set-TextValue 12 4 '0.574'
Set-TextValue 12 5 '  -6.79%  '
Set-TextValue 13 4 '47.06'
Set-TextValue 13 5 '  -2.85%  '
Set-TextValue 14 4 '0.0000273'
Set-TextValue 14 5 '  -3.56%  '
Set-TextValue 15 4 '4.074.76'
Set-TextValue 15 5 '  -3.35%  '
Set-TextValue 16 4 '8.35'
Set-TextValue 16 5 '  -7.14%  '
Set-TextValue 17 4 '608.84'
Set-TextValue 17 5 '  -9.47%  '
Set-TextValue 18 4 '69.049.86'
Set-TextValue 18 5 '  -2.42%  '
Set-TextValue 19 4 '3.504.50'
Set-TextValue 19 5 '  -3.57%  '
Set-TextValue 20 4 '0.120'
Set-TextValue 20 5 '  -2.35%  '
Set-TextValue 21 4 '17.33'
Set-TextValue 21 5 '  -2.68%  '
Set-TextValue 22 4 '11.01'
Set-TextValue 22 5 '  -4.30%  '
Set-TextValue 23 4 '0.879'
Set-TextValue 23 5 '  -6.86%  '
Set-TextValue 24 4 '15.52'
Set-TextValue 24 5 '  -9.66%  '
Set-TextValue 25 4 '96.22'
Set-TextValue 25 5 '  -3.72%  '
Set-TextValue 26 4 '3.81'
Set-TextValue 26 5 '  -2.92%  '
Set-TextValue 27 4 '1.00'
Set-TextValue 27 5 '  +0.02%  '
Set-TextValue 28 4 '2.59'
Set-TextValue 28 5 '  -7.14%  '
Set-TextValue 29 4 '9.09'
Set-TextValue 29 5 '  -7.44%  '
Set-TextValue 30 4 '32.42'
Set-TextValue 30 5 '  -6.17%  '
Set-TextValue 31 4 '8.46'
Set-TextValue 31 5 '  -7.41%  '
Set-TextValue 32 4 '3.09'
Set-TextValue 32 5 '  -5.29%  '
Set-TextValue 33 4 '1.31'
Set-TextValue 33 5 '  -6.03%  '
Set-TextValue 34 4 '6.85'
Set-TextValue 34 5 '  -9.39%  '
Set-TextValue 35 4 '623.37'
Set-TextValue 35 5 '  +8.82%  '
Set-TextValue 36 4 '10.67'
Set-TextValue 36 5 '  -3.60%  '
Set-TextValue 37 4 '3.48'
Set-TextValue 37 5 '  -12.04%  '
Set-TextValue 38 4 '0.101'
Set-TextValue 38 5 '  -5.73%  '
Set-TextValue 39 4 '56.86'
Set-TextValue 39 5 '  -2.74%  '
Set-TextValue 40 4 '1.00'
Set-TextValue 40 5 '  +0.05%  '
Set-TextValue 41 4 '0.0438'
Set-TextValue 41 5 '  -2.98%  '
Set-TextValue 42 4 '0.134'
Set-TextValue 42 5 '  -3.62%  '
Set-TextValue 43 4 '3.377.23'
Set-TextValue 43 5 '  -4.79%  '
Set-TextValue 44 4 '0.323'
Set-TextValue 44 5 '  -6.40%  '
Set-TextValue 45 4 '32.52'
Set-TextValue 45 5 '  -5.28%  '
Set-TextValue 46 4 '0.0₃0694'
Set-TextValue 46 5 '  -5.04%  '
Set-TextValue 47 4 '2.51'
Set-TextValue 47 5 '  -6.14%  '
Set-TextValue 48 4 '2.72'
Set-TextValue 48 5 '  -8.07%  '
Set-TextValue 49 4 '0.128'
Set-TextValue 49 5 '  -5.29%  '
Set-TextValue 50 4 '133.91'
Set-TextValue 50 5 '  -2.91%  '
Set-TextValue 51 4 '5.51'
Set-TextValue 51 5 '  +9.84%  '
